# Applies updated probability values to Sheet1 per commit:
# "added more games, sped up simulate game logic, and drafted optimization logic"
# This refreshes the computed transition-matrix percentages in rows 2-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.2016260162601626
    "C2" = 0.5398373983739837
    "J2" = 0.01300813008130081
    "P2" = 0.1414634146341463
    "S2" = 0.1040650406504065
    "B3" = 0.01146131805157593
    "C3" = 0.02865329512893983
    "J3" = 0.03724928366762178
    "P3" = 0.7048710601719198
    "S3" = 0.2177650429799427
    "J4" = 0.1046511627906977
    "P4" = 0.7674418604651163
    "S4" = 0.1279069767441861
    "J5" = 0.25
    "P5" = 0.5
    "S5" = 0.25
    "B6" = 0.07399577167019028
    "D6" = 0.0105708245243129
    "E6" = 0.004228329809725159
    "F6" = 0.05285412262156448
    "J6" = 0.2832980972515856
    "O6" = 0.03382663847780127
    "Q6" = 0.1585623678646934
    "R6" = 0.04651162790697674
    "S6" = 0.3361522198731501
    "B7" = 0.1129032258064516
    "D7" = 0.01382488479262673
    "F7" = 0.04838709677419355
    "J7" = 0.1797235023041475
    "O7" = 0.02073732718894009
    "Q7" = 0.1728110599078341
    "R7" = 0.05069124423963134
    "S7" = 0.4009216589861751
    "B8" = 0.09751243781094528
    "D8" = 0.01691542288557214
    "E8" = 0.0009950248756218905
    "F8" = 0.06368159203980099
    "J8" = 0.1194029850746269
    "O8" = 0.02288557213930348
    "Q8" = 0.1781094527363184
    "R8" = 0.06766169154228856
    "S8" = 0.4328358208955224
    "B9" = 0.1063829787234043
    "D9" = 0.0060790273556231
    "E9" = 0.00303951367781155
    "F9" = 0.0911854103343465
    "J9" = 0.121580547112462
    "O9" = 0.0182370820668693
    "Q9" = 0.1458966565349544
    "R9" = 0.0851063829787234
    "S9" = 0.4224924012158054
    "B10" = 0.1197806832560101
    "D10" = 0.02446225221425559
    "E10" = 0.0008435259384226065
    "F10" = 0.07043441585828764
    "J10" = 0.1037536904259806
    "O10" = 0.02446225221425559
    "Q10" = 0.2113032475748629
    "R10" = 0.06832560101223113
    "S10" = 0.3766343315056938
    "F11" = 0.001545595054095827
    "G11" = 0.160741885625966
    "J11" = 0.09428129829984544
    "K11" = 0.1993817619783617
    "L11" = 0.5239567233384853
    "S11" = 0.02009273570324575
    "G12" = 0.7442455242966752
    "J12" = 0.1662404092071611
    "K12" = 0.01534526854219949
    "L12" = 0.02813299232736573
    "S12" = 0.04603580562659847
    "G13" = 0.594059405940594
    "J13" = 0.297029702970297
    "S13" = 0.1089108910891089
    "F15" = 0.02428256070640177
    "H15" = 0.1920529801324503
    "I15" = 0.04194260485651214
    "J15" = 0.3134657836644592
    "K15" = 0.0772626931567329
    "M15" = 0.01324503311258278
    "N15" = 0.002207505518763797
    "O15" = 0.05960264900662252
    "S15" = 0.2759381898454746
    "F16" = 0.02813299232736573
    "H16" = 0.1457800511508952
    "I16" = 0.07928388746803069
    "J16" = 0.4092071611253197
    "K16" = 0.1355498721227621
    "M16" = 0.007672634271099744
    "O16" = 0.06393861892583121
    "S16" = 0.1304347826086956
    "F17" = 0.02655889145496536
    "H17" = 0.2032332563510393
    "I17" = 0.07159353348729793
    "J17" = 0.3937644341801386
    "K17" = 0.1039260969976905
    "M17" = 0.02309468822170901
    "N17" = 0.00115473441108545
    "O17" = 0.06351039260969978
    "S17" = 0.1131639722863741
    "F18" = 0.03642384105960265
    "H18" = 0.1754966887417219
    "I18" = 0.08940397350993377
    "J18" = 0.4105960264900662
    "K18" = 0.08278145695364239
    "M18" = 0.01986754966887417
    "O18" = 0.05629139072847682
    "S18" = 0.1291390728476821
    "F19" = 0.02821801314263626
    "H19" = 0.2365674526478547
    "I19" = 0.06957866254348666
    "J19" = 0.332044839582528
    "K19" = 0.1113258600695787
    "M19" = 0.02435253189022033
    "N19" = 0.001932740626207963
    "O19" = 0.06494008504058756
    "S19" = 0.1310398144568999
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
